$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Darshan"
$ws.Range("B3").Value = "Saman"
$ws.Range("B4").Value = "Shail"
$ws.Range("B5").Value = "Vaibhav"
$ws.Range("B6").Value = "John"
$ws.Range("B7").Value = "d"
$ws.Range("B8").Value = "e"
$ws.Range("B9").Value = "f"
$ws.Range("B10").Value = "g"
$ws.Range("B11").Value = "h"
$ws.Range("B12").Value = "i"

$ws.Range("C3").Value = "s.sh@husky.neu.edu"
$ws.Range("C4").Value = "s.sh@husky.neu.edu"
$ws.Range("C5").Value = "s.sh@husky.neu.edu"
$ws.Range("C6").Value = "s.sh@husky.neu.edu"
$ws.Range("C7").Value = "s.sh@husky.neu.edu"
$ws.Range("C8").Value = "s.sh@husky.neu.edu"
$ws.Range("C9").Value = "s.sh@husky.neu.edu"
$ws.Range("C10").Value = "s.sh@husky.neu.edu"
$ws.Range("C11").Value = "s.sh@husky.neu.edu"
$ws.Range("C12").Value = "s.sh@husky.neu.edu"

$ws.Range("C6").Select()
